$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Write the refreshed scrape data into rows 2-15 (rows 14-15 are brand new appended rows;
# two brand-new listings were inserted above older ones, shifting the remaining rows down).

# Row 2
$ws.Cells.Item(2,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(2,2).Value = "EC×AIプロダクト/業務改善リード"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5467702"
$ws.Cells.Item(2,7).Value = 338
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◇業務改善"

# Row 3
$ws.Cells.Item(3,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(3,2).Value = "初回 急募 自動カートインツール 開発のプロフェッショナルを探しています"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5467745"
$ws.Cells.Item(3,7).Value = 120
$ws.Cells.Item(3,8).Value = "◆ツール,開発"

# Row 4
$ws.Cells.Item(4,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(4,2).Value = "【緊急募集】動画解析アプリ開発のプロフェッショナル"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5467910"
$ws.Cells.Item(4,7).Value = 88
$ws.Cells.Item(4,8).Value = "◆開発 ◇アプリ"

# Row 5
$ws.Cells.Item(5,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(5,2).Value = "【急募】Ecommerce開発のシニアデベロッパーを探しています"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5467698"
$ws.Cells.Item(5,7).Value = 75
$ws.Cells.Item(5,8).Value = "◆開発"

# Row 6
$ws.Cells.Item(6,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(6,2).Value = "【法人歓迎】プローバステージ制御ソフト開発の見積依頼"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5467295"
$ws.Cells.Item(6,7).Value = 75
$ws.Cells.Item(6,8).Value = "◆開発"

# Row 7
$ws.Cells.Item(7,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(7,2).Value = "イベントサイトのWeb制作(決済機能付き)依頼"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5467460"
$ws.Cells.Item(7,7).Value = 38
$ws.Cells.Item(7,8).Value = "◇サイト"

# Row 8
$ws.Cells.Item(8,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(8,2).Value = "【急募】社内Webアプリの修正・再構築依頼"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5467384"
$ws.Cells.Item(8,7).Value = 33
$ws.Cells.Item(8,8).Value = "◇アプリ"

# Row 9
$ws.Cells.Item(9,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(9,2).Value = "iPhoneのブラウザ要素の書き換えアプリ作成"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5467578"
$ws.Cells.Item(9,7).Value = 30
$ws.Cells.Item(9,8).Value = "◇アプリ"

# Row 10
$ws.Cells.Item(10,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(10,2).Value = "進行管理およびチームディレクションを担当"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Cells.Item(10,7).Value = 30
$ws.Cells.Item(10,8).Value = "◇管理"

# Row 11
$ws.Cells.Item(11,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(11,2).Value = "限定公開 限定公開の仕事"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "5,000,000 円 ~ / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,6).Value = "https://www.lancers.jp/work/detail/5467882"
$ws.Cells.Item(11,7).Value = 25
$ws.Cells.Item(11,8).ClearContents()

# Row 12
$ws.Cells.Item(12,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(12,2).Value = "【急募】cloudflare導入の経験者を探しています!"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,6).Value = "https://www.lancers.jp/work/detail/5467334"
$ws.Cells.Item(12,7).Value = 18
$ws.Cells.Item(12,8).ClearContents()

# Row 13
$ws.Cells.Item(13,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(13,2).Value = "電気点火装置の回路図作成依頼"
$ws.Cells.Item(13,3).Value = "システム開発"
$ws.Cells.Item(13,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(13,5).Value = "期限情報なし"
$ws.Cells.Item(13,6).Value = "https://www.lancers.jp/work/detail/5466994"
$ws.Cells.Item(13,7).Value = 13
$ws.Cells.Item(13,8).ClearContents()

# Row 14
$ws.Cells.Item(14,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(14,2).Value = "リダイレクトでエラーが出てるので修正依頼"
$ws.Cells.Item(14,3).Value = "システム開発"
$ws.Cells.Item(14,4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(14,5).Value = "期限情報なし"
$ws.Cells.Item(14,6).Value = "https://www.lancers.jp/work/detail/5467706"
$ws.Cells.Item(14,7).Value = 10
$ws.Cells.Item(14,8).ClearContents()

# Row 15
$ws.Cells.Item(15,1).Value = "2026-01-09 02:01:12"
$ws.Cells.Item(15,2).Value = "ドメインの移行をして欲しい"
$ws.Cells.Item(15,3).Value = "システム開発"
$ws.Cells.Item(15,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(15,5).Value = "期限情報なし"
$ws.Cells.Item(15,6).Value = "https://www.lancers.jp/work/detail/5467598"
$ws.Cells.Item(15,7).Value = 10
$ws.Cells.Item(15,8).ClearContents()

# Register hyperlink relationships for the two newly appended rows, keeping all
# previously existing hyperlink entries (F2..F13) completely untouched.
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5467910")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5467882")
$ws.Cells.Item(14,6).Style = "Hyperlink"
$ws.Cells.Item(15,6).Style = "Hyperlink"

Write-Output "done"
